# Weekly update: insert two new "Repollo" price records (week of 2023-01-13)
# for "Feria Lagunitas de Puerto Montt" and shift the existing historical
# rows down by two positions (rows 599-623 -> 601-625).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 599, pushing the
# existing data (and the two trailing rows that fall off the old R623
# dimension) down to rows 601-625.
$ws.Range("A599:A600").EntireRow.Insert()

# --- New row 599 ---
$ws.Cells.Item(599, 1).Value = 4
$ws.Cells.Item(599, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(599, 3).Value = "Los Lagos"
$ws.Cells.Item(599, 4).Value = 44939
$ws.Cells.Item(599, 5).Value = 10
$ws.Cells.Item(599, 6).Value = 100112006
$ws.Cells.Item(599, 7).Value = "Repollo"
$ws.Cells.Item(599, 8).Value = "Copenhague"
$ws.Cells.Item(599, 9).Value = "Primera"
$ws.Cells.Item(599, 10).Value = 600
$ws.Cells.Item(599, 11).Value = 2000
$ws.Cells.Item(599, 12).Value = 2000
$ws.Cells.Item(599, 13).Value = 2000
$ws.Cells.Item(599, 14).Value = "$/unidad"
$ws.Cells.Item(599, 15).Value = "Región Metropolitana"
$ws.Cells.Item(599, 16).Value = 2000
$ws.Cells.Item(599, 17).Value = 1
$ws.Cells.Item(599, 18).Value = "Hortaliza"

# --- New row 600 ---
$ws.Cells.Item(600, 1).Value = 4
$ws.Cells.Item(600, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(600, 3).Value = "Los Lagos"
$ws.Cells.Item(600, 4).Value = 44939
$ws.Cells.Item(600, 5).Value = 10
$ws.Cells.Item(600, 6).Value = 100112006
$ws.Cells.Item(600, 7).Value = "Repollo"
$ws.Cells.Item(600, 8).Value = "Crespo record"
$ws.Cells.Item(600, 9).Value = "Primera"
$ws.Cells.Item(600, 10).Value = 600
$ws.Cells.Item(600, 11).Value = 1800
$ws.Cells.Item(600, 12).Value = 1800
$ws.Cells.Item(600, 13).Value = 1800
$ws.Cells.Item(600, 14).Value = "$/unidad"
$ws.Cells.Item(600, 15).Value = "Región Metropolitana"
$ws.Cells.Item(600, 16).Value = 1800
$ws.Cells.Item(600, 17).Value = 1
$ws.Cells.Item(600, 18).Value = "Hortaliza"
